# "remove part about sue" -- drop the redundant/duplicated passage that
# followed "...I'll take care of you." (the bit about doors held open,
# cleaning tables, and keeping five dollars in the car), since the same
# sentiment is already covered by the "In their love that they preached
# and lived..." paragraph that immediately follows it.

$d = $word.ActiveDocument

$rightQuote = [char]0x2019

$old = " I didn" + $rightQuote + "t just see how they loved me, but how they loved the" `
     + " people around them: in the doors they held open and in the effort they put into" `
     + " cleaning the table before leaving the restaurant. Another sister told me how she" `
     + " always kept five dollars in cash in her car, so she would always have something to give."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

if (-not $found) {
    throw "Could not find the passage to remove."
}

Write-Output "Removed the duplicate passage after `"...take care of you.`""
